# Rewrite the Sheet1 table with the new cards data (header unchanged),
# growing the table from 4 rows to 8 rows, and move the selection to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (row 1 = header, rows 2-8 = data).
# Row 8's Description/Tags are plain numbers (no t="s"), matching the source edit.
$rows = @(
    @('Name',      'Description', 'Tags'),
    @('gkj;t43',   'hsdgh',       'sh'),
    @('Gew3',      '2tg',         'ghf'),
    @('RWQ`',      'afgd',        'sgsfsf'),
    @('FSDg',      'gdfeR',       'sfgs'),
    @('GDF',       'G',           'fgaf'),
    @('DA',        'FGADF',       'dg3'),
    @('GAGADF',    124,           4234)
)

for ($r = 1; $r -le $rows.Length; $r++) {
    $rowData = $rows[$r - 1]
    for ($c = 1; $c -le 3; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Move the active selection to A4, matching the post-edit cursor position.
[void]$ws.Range("A4").Select()
